# Test_Specification_Template.xlsx edit
# Commit: Removed select_Itam function as the click_button function has
# been updated to support this.
#
# The "Lists" sheet holds the list of supported Ruby_Web_Functions in
# column C (with optional notes in column D). Remove the rows for the
# functions that were retired: Open_Portal_URL, Select_Item, Portal_Login,
# SINT_Login and Admin_Portal_Login. Deleting the whole rows (bottom-to-top
# so earlier row numbers stay valid) shifts everything below up, which is
# exactly what the target workbook shows (rows 3-29 instead of 3-34).

$wb = $excel.ActiveWorkbook
$lists = $wb.Worksheets.Item("Lists")

$lists.Rows.Item(34).Delete()
$lists.Rows.Item(33).Delete()
$lists.Rows.Item(32).Delete()
$lists.Rows.Item(27).Delete()
$lists.Rows.Item(22).Delete()

# Leave the cursor roughly where the author left it (Lists!B37, Sheet1!C5)
$lists.Range("B37").Select()
$wb.Worksheets.Item("Sheet1").Range("C5").Select()

# The author's last active sheet was "Lists"
$lists.Activate()
